# Updates cryptos list values to match the latest scrape (GitHub Actions run).
# Column D ("Price") cells that parse as plain numbers are written with a leading
# apostrophe so Excel stores them as text (matching the sheet's existing inlineStr
# convention) instead of silently converting them to numeric values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


$ws.Range("D2").Value = "60.852.36"
$ws.Range("E2").Value = "  +0.62%  "

$ws.Range("D3").Value = "2.593.42"
$ws.Range("E3").Value = "  +0.33%  "

$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").Value = "'522.39"
$ws.Range("E5").Value = "  +2.98%  "

$ws.Range("D6").Value = "'154.10"
$ws.Range("E6").Value = "  +0.43%  "

$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("D8").Value = "'0.595"
$ws.Range("E8").Value = "  +2.75%  "

$ws.Range("E9").Value = "  +2.03%  "

$ws.Range("E10").Value = "  +1.34%  "

$ws.Range("E11").Value = "  -0.31%  "

$ws.Range("E12").Value = "  +1.39%  "

$ws.Range("D13").Value = "3.048.74"
$ws.Range("E13").Value = "  +0.26%  "

$ws.Range("D14").Value = "60.876.20"
$ws.Range("E14").Value = "  +0.66%  "

$ws.Range("D15").Value = "'21.66"
$ws.Range("E15").Value = "  +0.22%  "

$ws.Range("E16").Value = "  +0.32%  "

$ws.Range("D17").Value = "2.601.12"
$ws.Range("E17").Value = "  +0.25%  "

$ws.Range("D18").Value = "'4.76"
$ws.Range("E18").Value = "  -0.51%  "

$ws.Range("D19").Value = "'353.66"
$ws.Range("E19").Value = "  +2.34%  "

$ws.Range("E20").Value = "  +1.34%  "

$ws.Range("E21").Value = "  +1.43%  "

$ws.Range("E22").Value = "  +0.28%  "

$ws.Range("D23").Value = "'60.81"
$ws.Range("E23").Value = "  +1.44%  "

$ws.Range("E24").Value = "  +1.63%  "

$ws.Range("E25").Value = "  -0.26%  "

$ws.Range("D26").Value = "2.708.16"
$ws.Range("E26").Value = "  +0.23%  "

$ws.Range("D27").Value = "'0.998"
$ws.Range("E27").Value = "  -0.09%  "

$ws.Range("D28").Value = "0.0₃0845"
$ws.Range("E28").Value = "  +0.30%  "

$ws.Range("D29").Value = "'7.38"
$ws.Range("E29").Value = "  +0.08%  "

$ws.Range("E30").Value = "  -0.07%  "

$ws.Range("D31").Value = "'6.35"
$ws.Range("E31").Value = "  +11.14%  "

$ws.Range("D32").Value = "'19.37"
$ws.Range("E32").Value = "  +0.15%  "

$ws.Range("E33").Value = "  +2.63%  "

$ws.Range("D34").Value = "'148.07"
$ws.Range("E34").Value = "  -3.67%  "

$ws.Range("E35").Value = "  +4.28%  "

$ws.Range("D36").Value = "'0.932"
$ws.Range("E36").Value = "  +9.18%  "

$ws.Range("E37").Value = "  +0.96%  "

$ws.Range("D38").Value = "'0.861"
$ws.Range("E38").Value = "  +1.33%  "

$ws.Range("E39").Value = "  +1.87%  "

$ws.Range("E40").Value = "  +1.52%  "

$ws.Range("E41").Value = "  +1.40%  "

$ws.Range("D42").Value = "'288.38"
$ws.Range("E42").Value = "  -2.21%  "

$ws.Range("D43").Value = "'0.101"
$ws.Range("E43").Value = "  +1.98%  "

$ws.Range("D44").Value = "'0.621"
$ws.Range("E44").Value = "  -0.55%  "

$ws.Range("D45").Value = "'0.0560"
$ws.Range("E45").Value = "  +0.68%  "

$ws.Range("D46").Value = "'0.998"
$ws.Range("E46").Value = "  +0.01%  "

$ws.Range("E47").Value = "  -1.20%  "

$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "'4.88"
$ws.Range("E48").Value = "  +0.60%  "

$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").Value = "'0.0238"
$ws.Range("E49").Value = "  +2.03%  "

$ws.Range("E50").Value = "  +0.17%  "

$ws.Range("D51").Value = "'19.11"
$ws.Range("E51").Value = "  +8.23%  "
